$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 54, pushing the existing rows 54-149
# down to 56-151 (this matches Excel's native row-insert behaviour and keeps
# formatting/styles of column D, which carries a date number format).
$ws.Rows("54:55").Insert()

# Populate the two newly inserted rows (54 and 55) with the new price records.
# Columns A, B, C, E, F, G, N, Q, R hold the same constant values used
# throughout the rest of the sheet (market/region/category/unit metadata).

# Row 54
$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value = "Bíobío"
$ws.Cells.Item(54, 4).Value = 44540
$ws.Cells.Item(54, 5).Value = 8
$ws.Cells.Item(54, 6).Value = 100114001
$ws.Cells.Item(54, 7).Value = "Papa"
$ws.Cells.Item(54, 8).Value = "Asterix"
$ws.Cells.Item(54, 9).Value = "1a nueva(o)"
$ws.Cells.Item(54, 10).Value = 250
$ws.Cells.Item(54, 11).Value = 9500
$ws.Cells.Item(54, 12).Value = 10000
$ws.Cells.Item(54, 13).Value = 9700
$ws.Cells.Item(54, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(54, 16).Value = 388
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"

# Row 55
$ws.Cells.Item(55, 1).Value = 11
$ws.Cells.Item(55, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value = "Bíobío"
$ws.Cells.Item(55, 4).Value = 44540
$ws.Cells.Item(55, 5).Value = 8
$ws.Cells.Item(55, 6).Value = 100114001
$ws.Cells.Item(55, 7).Value = "Papa"
$ws.Cells.Item(55, 8).Value = "Patagonia"
$ws.Cells.Item(55, 9).Value = "1a nueva(o)"
$ws.Cells.Item(55, 10).Value = 270
$ws.Cells.Item(55, 11).Value = 8000
$ws.Cells.Item(55, 12).Value = 9500
$ws.Cells.Item(55, 13).Value = 8667
$ws.Cells.Item(55, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(55, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(55, 16).Value = 347
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = "Hortaliza"
